$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.452.72'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.223.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -6.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.63'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.20%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.37'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -8.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.06'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.59'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -9.00%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.26%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.86'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.556.72'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.231.89'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.352.15'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.15'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.38%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.04'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.38'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +13.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.92'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.94'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.38'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -7.13%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.91'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.09'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -11.39%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.15'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.80'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0873'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.57'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.126'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.17'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.75'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.229'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.56%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -11.45%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.31'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -7.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +10.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.27'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.06'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.40'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.20%  '
